$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values must remain plain text, matching the
# original inline-string cells in the sheet, even though several look numeric)
$updates = [ordered]@{
    'D2' = '29.908.48'
    'E2' = '  -0.15%  '
    'D3' = '1.875.35'
    'E3' = '  -0.95%  '
    'D4' = '0.9995'
    'E4' = '  -0.07%  '
    'D5' = '0.7395'
    'E5' = '  -5.03%  '
    'D6' = '242.60'
    'E6' = '  -0.54%  '
    'D7' = '0.9995'
    'E7' = '  -0.09%  '
    'D8' = '0.3164'
    'E8' = '  +0.89%  '
    'D9' = '0.07213'
    'E9' = '  -0.88%  '
    'D10' = '24.77'
    'E10' = '  -4.19%  '
    'D11' = '0.08400'
    'E11' = '  -3.28%  '
    'D12' = '0.7506'
    'E12' = '  -3.06%  '
    'D13' = '5.427'
    'E13' = '  +0.33%  '
    'D14' = '1.887.28'
    'E14' = '  -13.14%  '
    'D15' = '92.57'
    'E15' = '  -2.09%  '
    'D16' = '29.904.77'
    'E16' = '  -0.69%  '
    'E17' = '  -1.86%  '
    'E18' = '  -2.60%  '
    'D19' = '243.62'
    'E19' = '  -0.91%  '
    'E20' = '  -0.64%  '
    'D21' = '0.9994'
    'E21' = '  -0.14%  '
    'D22' = '2.123.27'
    'E22' = '  -6.82%  '
    'E23' = '  -2.82%  '
    'D24' = '1.004'
    'E24' = '  +0.31%  '
    'E25' = '  -6.83%  '
    'D26' = '9.278'
    'E26' = '  -2.36%  '
    'D27' = '165.38'
    'E27' = '  +1.23%  '
    'E28' = '  -1.30%  '
    'D29' = '2.033'
    'E29' = '  -0.97%  '
    'D30' = '1.508'
    'E30' = '  +5.24%  '
    'D31' = '4.592'
    'E31' = '  +1.77%  '
    'D32' = '1.532'
    'E32' = '  -0.68%  '
    'D33' = '4.268'
    'E33' = '  +3.33%  '
    'E34' = '  -3.26%  '
    'E35' = '  -0.93%  '
    'D36' = '0.7543'
    'E36' = '  +0.00%  '
    'D37' = '0.9974'
    'E37' = '  -0.44%  '
    'D38' = '2.701'
    'E38' = '  +0.58%  '
    'D39' = '0.01960'
    'E39' = '  -0.08%  '
    'D40' = '2.754'
    'E40' = '  -1.40%  '
    'D41' = '0.4534'
    'E41' = '  +0.44%  '
    'D42' = '1.117.36'
    'E42' = '  +0.76%  '
    'D43' = '6.035'
    'E43' = '  -0.68%  '
    'D44' = '72.52'
    'E44' = '  -1.62%  '
    'D45' = '0.8571'
    'E45' = '  +0.46%  '
    'E46' = '  +0.08%  '
    'E47' = '  +0.15%  '
    'D48' = '3.122'
    'E48' = '  +3.82%  '
    'D49' = '7.633'
    'E49' = '  +0.40%  '
    'D50' = '1.840'
    'E50' = '  -2.20%  '
    'D51' = '2.020.44'
    'E51' = '  -6.36%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
